$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Beveilingseisen" (security requirements) column (E) with the
# new requirement texts. The three brand-new strings are introduced in the
# order Cross-Site Request Forgery -> Authentication -> Broken Access
# Control (first touched on E10, E2, E3 respectively) so the rebuilt shared
# string table lands in the same order as the target workbook.
$ws.Range("E10").Value = "Cross-Site Request Forgery"
$ws.Range("E2").Value = "Authentication"
$ws.Range("E3").Value = "Broken Access Control"
$ws.Range("E4").Value = "Authentication"
$ws.Range("E5").Value = "Tegen SQL en XSS beveiligd"
$ws.Range("E6").Value = "Geen"
$ws.Range("E7").Value = "Geen"
$ws.Range("E8").Value = "Geen"
$ws.Range("E9").Value = "Geen"
$ws.Range("E11").Value = "Geen"
$ws.Range("E12").Value = "Geen"
$ws.Range("E13").Value = "Geen"

# Move the sheet's active selection from F4 to D13.
$ws.Range("D13").Select()
